$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 105 - OpenVibSpec
$ws.Range("A105").Value = "OpenVibSpec"
$ws.Range("B105").Value = "https://github.com/arnrau/OpenVibSpec"
$ws.Range("C105").Value = "https://github.com/arnrau/OpenVibSpec"
$ws.Range("E105").Value = "Python"
$ws.Range("F105").Value = "Vibrational spectroscopy"
$ws.Range("G105").Value = "IR, Raman"

# Row 106 - octavvs
$ws.Range("A106").Value = "octavvs"
$ws.Range("B106").Value = "https://github.com/ctroein/octavvs"
$ws.Range("C106").Value = "https://github.com/ctroein/octavvs"
$ws.Range("E106").Value = "Python"
$ws.Range("F106").Value = "Vibrational spectroscopy"
$ws.Range("G106").Value = "IR, Raman"

# Mirror the author's final selection position (next empty row below the new data)
[void]$ws.Range("A107").Select()
